$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new column H (copy style from G1 header so it matches the other headers)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Label"

# Set Label column values: first 5 patients (Control) = 0, next 5 (MDD) = 1.
# This pattern repeats for each block of patient rows (100-iteration block
# starting at row 2, 200-iteration block starting at row 12).
$labels = @(0,0,0,0,0,1,1,1,1,1)
foreach ($blockStart in @(2,12)) {
    for ($i = 0; $i -lt 10; $i++) {
        $row = $blockStart + $i
        $ws.Cells.Item($row, 8).Value = $labels[$i]
    }
}
